# Apply "Document updates and additions" commit:
#  - Update status on the "9-30-13" sheet (5th sheet / sheet5.xml):
#      * Row 4 ("Android to Vex Comm"): record a completion date in column C
#        and bump % complete to 100%.
#      * Row 5 ("Sockets Tutorial"): record a meeting date in column M
#        (same date format as the existing M4 entry).
#      * Row 6 ("Interrupt based obstacle avoidance"): record a completion
#        date in column C, bump % complete to 100% and update hours spent.
#      * Add a new task row 7, "Project Plan 1" - the Phase 2 plan task,
#        with its due date, % complete and hours spent.
#  - Leave the active selection on N5, matching the last cell touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)   # "9-30-13" -> xl/worksheets/sheet5.xml

# --- Row 4: Android to Vex Comm -------------------------------------------
$ws.Range("C4").Value = 40075      # completion date (1904 date system serial)
$ws.Range("D4").Value = 1          # % complete -> 100%

# --- Row 5: Sockets Tutorial -------------------------------------------
# Copy the date formatting already used for M4 so the new M5 cell shares
# the same style instead of creating a brand new number format.
$ws.Range("M4").Copy()
$ws.Range("M5").PasteSpecial(-4122)
$ws.Range("M5").Value = 40085

# --- Row 6: Interrupt based obstacle avoidance -----------------------------
$ws.Range("C6").Value = 40080      # completion date
$ws.Range("D6").Value = 1          # % complete -> 100%
$ws.Range("E6").Value = 5.5        # hours spent

# --- Row 7 (new): Project Plan 1 -------------------------------------------
$ws.Range("A7").Value = "Project Plan 1"
$ws.Range("B7").Value = 40086      # due date
$ws.Range("D7").Value = 0.9        # % complete -> 90%
$ws.Range("E7").Value = 3.5        # hours spent

# --- Update active selection to reflect where work left off ---------------
$ws.Activate()
$ws.Range("N5").Select()
